$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 6-23 in the original sheet held card codes 26KBF1KVM005..26KBF1KVM022
# with only column A (and, for row 6, B/C) populated -- row 5 was blank.
# This revision inserts the missing row 5 (card 26KBF1KVM005 / "Maalamma"),
# shifts every subsequent record up by one row, and backfills the
# previously-empty amount/status/type columns (C/D/E) for every record.
# The final record (26KBF1KVM022) also gets a "mpg" name and an amount of
# 4000 instead of 2000.

$ws.Range("A5").Value = "26KBF1KVM005"
$ws.Range("B5").Value = "Maalamma"
$ws.Range("C5").Value = 2000
$ws.Range("D5").Value = "CLEARED"
$ws.Range("E5").Value = "KVM"
$ws.Range("A6").Value = "26KBF1KVM006"
$ws.Range("B6").Value = "Lakshmi Mahadevamma 1"
$ws.Range("C6").Value = 2000
$ws.Range("D6").Value = "CLEARED"
$ws.Range("E6").Value = "KVM"
$ws.Range("A7").Value = "26KBF1KVM007"
$ws.Range("B7").Value = "Manjanna 1"
$ws.Range("C7").Value = 2000
$ws.Range("D7").Value = "CLEARED"
$ws.Range("E7").Value = "KVM"
$ws.Range("A8").Value = "26KBF1KVM008"
$ws.Range("B8").Value = "Manjanna 2"
$ws.Range("C8").Value = 2000
$ws.Range("D8").Value = "CLEARED"
$ws.Range("E8").Value = "KVM"
$ws.Range("A9").Value = "26KBF1KVM009"
$ws.Range("B9").Value = "Manjanna 3"
$ws.Range("C9").Value = 2000
$ws.Range("D9").Value = "CLEARED"
$ws.Range("E9").Value = "KVM"
$ws.Range("A10").Value = "26KBF1KVM010"
$ws.Range("B10").Value = "Shivanna"
$ws.Range("C10").Value = 2000
$ws.Range("D10").Value = "CLEARED"
$ws.Range("E10").Value = "KVM"
$ws.Range("A11").Value = "26KBF1KVM011"
$ws.Range("B11").Value = "Srinivas leelamma 1"
$ws.Range("C11").Value = 2000
$ws.Range("D11").Value = "CLEARED"
$ws.Range("E11").Value = "KVM"
$ws.Range("A12").Value = "26KBF1KVM012"
$ws.Range("B12").Value = "Srinivas leelamma 2"
$ws.Range("C12").Value = 2000
$ws.Range("D12").Value = "CLEARED"
$ws.Range("E12").Value = "KVM"
$ws.Range("A13").Value = "26KBF1KVM013"
$ws.Range("B13").Value = "Srinivas leelamma 3"
$ws.Range("C13").Value = 2000
$ws.Range("D13").Value = "CLEARED"
$ws.Range("E13").Value = "KVM"
$ws.Range("A14").Value = "26KBF1KVM014"
$ws.Range("B14").Value = "Srinivas leelamma 4"
$ws.Range("C14").Value = 2000
$ws.Range("D14").Value = "CLEARED"
$ws.Range("E14").Value = "KVM"
$ws.Range("A15").Value = "26KBF1KVM015"
$ws.Range("B15").Value = "Srinivas leelamma 5"
$ws.Range("C15").Value = 2000
$ws.Range("D15").Value = "CLEARED"
$ws.Range("E15").Value = "KVM"
$ws.Range("A16").Value = "26KBF1KVM016"
$ws.Range("B16").Value = "Suma Mahadevamma 1"
$ws.Range("C16").Value = 2000
$ws.Range("D16").Value = "CLEARED"
$ws.Range("E16").Value = "KVM"
$ws.Range("A17").Value = "26KBF1KVM017"
$ws.Range("B17").Value = "Suma Mahadevamma 2"
$ws.Range("C17").Value = 2000
$ws.Range("D17").Value = "CLEARED"
$ws.Range("E17").Value = "KVM"
$ws.Range("A18").Value = "26KBF1KVM018"
$ws.Range("B18").Value = "Suma Muddanayaka"
$ws.Range("C18").Value = 2000
$ws.Range("D18").Value = "CLEARED"
$ws.Range("E18").Value = "KVM"
$ws.Range("A19").Value = "26KBF1KVM019"
$ws.Range("B19").Value = "Muddanayaka Mahadevi"
$ws.Range("C19").Value = 2000
$ws.Range("D19").Value = "CLEARED"
$ws.Range("E19").Value = "KVM"
$ws.Range("A20").Value = "26KBF1KVM020"
$ws.Range("B20").Value = "Bhagya "
$ws.Range("C20").Value = 2000
$ws.Range("D20").Value = "CLEARED"
$ws.Range("E20").Value = "KVM"
$ws.Range("A21").Value = "26KBF1KVM021"
$ws.Range("B21").Value = "mpg"
$ws.Range("C21").Value = 2000
$ws.Range("D21").Value = "CLEARED"
$ws.Range("E21").Value = "KVM"
$ws.Range("A22").Value = "26KBF1KVM022"
$ws.Range("B22").Value = "mpg"
$ws.Range("C22").Value = 4000
$ws.Range("D22").Value = "CLEARED"
$ws.Range("E22").Value = "KVM"

# Row 23 no longer exists in the rebuilt table (data now ends at row 22).
$ws.Rows.Item(23).ClearContents()

# Match the saved selection/active cell from the edited workbook.
[void]$ws.Range("I16").Select()
